$d = $word.ActiveDocument

# 1) Update the revision date/time in the Date-styled paragraph.
$d.Content.Find.Execute(
    "June  16, 2021 (08:20:31 PM)", $true, $false, $false, $false, $false,
    $true, 1, $false, "June  16, 2021 (08:22:28 PM)", 2) | Out-Null

# 2) Split the "After completing the table..." sentence into several runs,
#    quoting "returns" and "value" with curly double quotes, while keeping
#    the paragraph's existing style (Definition).
$rng = $d.Content
$rng.Find.Execute(
    "After completing the table, can you detect a pattern between return type and resulting value?") | Out-Null

$lquote = [char]0x201C
$rquote = [char]0x201D

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:pPr><w:pStyle w:val="Definition"/></w:pPr>' +
       '<w:r><w:t xml:space="preserve">After completing the table, can you detect a pattern between</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">' + $lquote + '</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">returns</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">' + $rquote + '</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">and</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">' + $lquote + '</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">value</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">' + $rquote + '</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">?</w:t></w:r>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$rng.Collapse(1)
$rng.InsertXML($xml)
